$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 583, shifting the existing rows
# (old 583-609) down to (585-611).
$ws.Rows("583:584").Insert()

# New row 583: Betarraga, Primera, weekly price update.
$ws.Range("A583").Value = 7
$ws.Range("B583").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C583").Value = "Ñuble"
$ws.Range("D583").Value = 45041
$ws.Range("E583").Value = 16
$ws.Range("F583").Value = 100114014
$ws.Range("G583").Value = "Betarraga"
$ws.Range("H583").Value = "Sin especificar"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 500
$ws.Range("K583").Value = 1000
$ws.Range("L583").Value = 1000
$ws.Range("M583").Value = 1000
$ws.Range("N583").Value = "$/paquete 5 unidades"
$ws.Range("O583").Value = "Provincia de Diguillín"
$ws.Range("P583").Value = 200
$ws.Range("Q583").Value = 5
$ws.Range("R583").Value = "Hortaliza"

# New row 584: Betarraga, Segunda, weekly price update.
$ws.Range("A584").Value = 7
$ws.Range("B584").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C584").Value = "Ñuble"
$ws.Range("D584").Value = 45041
$ws.Range("E584").Value = 16
$ws.Range("F584").Value = 100114014
$ws.Range("G584").Value = "Betarraga"
$ws.Range("H584").Value = "Sin especificar"
$ws.Range("I584").Value = "Segunda"
$ws.Range("J584").Value = 500
$ws.Range("K584").Value = 800
$ws.Range("L584").Value = 800
$ws.Range("M584").Value = 800
$ws.Range("N584").Value = "$/paquete 5 unidades"
$ws.Range("O584").Value = "Provincia de Diguillín"
$ws.Range("P584").Value = 160
$ws.Range("Q584").Value = 5
$ws.Range("R584").Value = "Hortaliza"
